# Commit 456a3b4 adds a newly-scraped Bilibili-show listing — "丽水·幻梦动漫嘉年华"
# (2024-11-30) — as row 2 on the sheets that track this exhibition-type
# event: 展览 (Exhibition) and 全部类型 (All types). Both sheets gain the
# identical row and their used range grows from A1:I1 to A1:I2.

$wb = $excel.ActiveWorkbook

$targetSheets = @("展览", "全部类型")

foreach ($sheetName in $targetSheets) {
    $ws = $wb.Worksheets.Item($sheetName)

    # A2: serial number, numeric, formatted like the header's A1 cell
    # (bold / bordered / centered style) — copy that formatting over.
    $ws.Cells.Item(1, 1).Copy()
    $ws.Cells.Item(2, 1).PasteSpecial(-4122)  # xlPasteFormats
    $ws.Cells.Item(2, 1).Value = 1

    # B2 looks like a date ("2024-11-30"); force it to stay plain text
    # (matching the source data, which stores it as a literal string) and
    # then drop back to the default "Normal" style so no stray number
    # format sticks to the cell.
    $ws.Cells.Item(2, 2).NumberFormat = "@"
    $ws.Cells.Item(2, 2).Value = "2024-11-30"
    $ws.Cells.Item(2, 2).Style = "Normal"

    $ws.Cells.Item(2, 3).Value = "丽水·幻梦动漫嘉年华"
    $ws.Cells.Item(2, 4).Value = "中东路848号(解放街交汇) 飞达国际大酒店"
    $ws.Cells.Item(2, 5).Value = "2024.11.30 09:00-11.30 16:30"
    $ws.Cells.Item(2, 6).Value = 0
    $ws.Cells.Item(2, 7).Value = 50
    $ws.Cells.Item(2, 8).Value = "https://show.bilibili.com/platform/detail.html?id=93730"
    $ws.Cells.Item(2, 9).Value = "//i0.hdslb.com/bfs/openplatform/202410/FCA7Mt5G1729319544606.jpeg"
}
